# Update dfc-gas-trip-v1.xlsx "survey" sheet to use a "likert" appearance
# for the satisfaction question (row 7):
#   - Add a new "appearance" column (G) with header in G1.
#   - Move the existing hint text ("Please indicate the extent to which
#     you agree or disagree.") from D7 (hint) into the new G7 (appearance).
#   - Replace D7 (hint) with "likert".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# New "appearance" column header
$ws.Range("G1").Value = "appearance"

# Preserve the existing hint text, then move it to the new appearance column
$existingHint = $ws.Range("D7").Text

$ws.Range("G7").Value = $existingHint
$ws.Range("D7").Value = "likert"
